# Insert a new data row at row 87 (pushing existing rows 87..174 down to 88..175)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above the current row 87. This shifts formatting down
# along with the cell contents (column D keeps its date style), matching the
# behaviour of Excel's own "Insert Sheet Rows".
$ws.Rows("87:87").Insert()

# Populate the newly inserted row with the new record values.
$ws.Cells.Item(87, 1).Value = 11
$ws.Cells.Item(87, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(87, 3).Value = "Bíobío"
$ws.Cells.Item(87, 4).Value = 44904
$ws.Cells.Item(87, 5).Value = 8
$ws.Cells.Item(87, 6).Value = 100112032
$ws.Cells.Item(87, 7).Value = "Zapallo italiano"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 450
$ws.Cells.Item(87, 11).Value = 6000
$ws.Cells.Item(87, 12).Value = 6500
$ws.Cells.Item(87, 13).Value = 6278
$ws.Cells.Item(87, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(87, 15).Value = "Región Metropolitana"
$ws.Cells.Item(87, 16).Value = 126
$ws.Cells.Item(87, 17).Value = 50
$ws.Cells.Item(87, 18).Value = "Hortaliza"
